# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting refreshed output data (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Cells.Item(5, 6).Value = 6350
$wsExhibit.Cells.Item(10, 6).Value = 315
$wsExhibit.Cells.Item(12, 6).Value = 678
$wsExhibit.Cells.Item(16, 6).Value = 405
$wsExhibit.Cells.Item(17, 6).Value = 48
$wsExhibit.Cells.Item(21, 6).Value = 375
$wsExhibit.Cells.Item(24, 6).Value = 1064
$wsExhibit.Cells.Item(25, 6).Value = 126
$wsExhibit.Cells.Item(31, 6).Value = 3533

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(20, 6).Value = 4084
$wsShow.Cells.Item(32, 6).Value = 1599

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(11, 6).Value = 6350
$wsAll.Cells.Item(19, 6).Value = 678
$wsAll.Cells.Item(25, 6).Value = 405
$wsAll.Cells.Item(28, 6).Value = 48
$wsAll.Cells.Item(34, 6).Value = 375
$wsAll.Cells.Item(44, 6).Value = 1599
$wsAll.Cells.Item(49, 6).Value = 3533
